# Edit script: insert 3 new weekly price rows for "Sin especificar" Plátano
# (dated 2021-09-27 / serial 44466) right before the existing block that used
# to start at row 240, pushing the rest of the table down by 3 rows.
# This mirrors the commit "Fruta / hortaliza, semanal" which adds a new
# week's worth of data (3 quality rows: Pintón, Primera Maduro, Primera Pintón).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three blank rows at 240, shifting old rows 240:343 down to 243:346.
$ws.Rows("240:242").Insert()

# 2) The rows that used to be 240:242 are now at 243:245. Duplicate their
#    content up into the new rows 240:242 so every column (A..T) is populated
#    with the correct static data (market, product, variety, unit, origin...).
$ws.Range("A243:T243").Copy($ws.Range("A240:T240"))
$ws.Range("A244:T244").Copy($ws.Range("A241:T241"))
$ws.Range("A245:T245").Copy($ws.Range("A242:T242"))

# 3) Overwrite the cells that actually differ for the new week: Fecha (D),
#    Precio mínimo (N), Precio máximo (O), Precio promedio ponderado (P) and
#    Precio $/Kg (S).
$ws.Cells.Item(240, 4).Value2 = 44466
$ws.Cells.Item(240, 14).Value2 = 15000
$ws.Cells.Item(240, 15).Value2 = 15000
$ws.Cells.Item(240, 16).Value2 = 15000
$ws.Cells.Item(240, 19).Value2 = 750

$ws.Cells.Item(241, 4).Value2 = 44466
$ws.Cells.Item(241, 14).Value2 = 16500
$ws.Cells.Item(241, 15).Value2 = 16500
$ws.Cells.Item(241, 16).Value2 = 16500
$ws.Cells.Item(241, 19).Value2 = 825

$ws.Cells.Item(242, 4).Value2 = 44466
$ws.Cells.Item(242, 14).Value2 = 17000
$ws.Cells.Item(242, 15).Value2 = 17000
$ws.Cells.Item(242, 16).Value2 = 17000
$ws.Cells.Item(242, 19).Value2 = 850
